$wb = $excel.ActiveWorkbook

# --- Sheet "string" (sheet1): varname, idvalue, value, valuecurrent, initials, notes
#     -> make, id, varname, value, valuecurrent, initials, notes
$ws1 = $wb.Worksheets.Item("string")
$ws1.Columns("B:B").Delete()
$ws1.Columns("A:B").Insert()
$ws1.Range("A1").Value = "make"
$ws1.Range("B1").Value = "id"
$ws1.Range("A2:B2").NumberFormat = "0"

# --- Sheet "numeric" (sheet2): varname, idvalue, value, valuecurrent, initials, notes
#     -> make, id, varname, value, valuecurrent, initials, notes
$ws2 = $wb.Worksheets.Item("numeric")
$ws2.Columns("B:B").Delete()
$ws2.Columns("A:B").Insert()
$ws2.Range("A1").Value = "make"
$ws2.Range("B1").Value = "id"
$ws2.Range("A2:B2").NumberFormat = "0"

# --- Sheet "drop" (sheet3): idvalue, initials, notes
#     -> make, id, n_obs, initials, notes
$ws3 = $wb.Worksheets.Item("drop")
$ws3.Columns("A:A").Delete()
$ws3.Columns("A:C").Insert()
$ws3.Range("A1").Value = "make"
$ws3.Range("B1").Value = "id"
$ws3.Range("C1").Value = "n_obs"
$ws3.Range("A2:C2").NumberFormat = "0"
